$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K8").Value = 12.67039049919483
$ws.Range("K9").Value = 12.67039049919483
$ws.Range("R9").Value = 3.135414186445925
$ws.Range("S9").Value = 3.817068050129936
$ws.Range("K10").Value = 12.67039049919483
$ws.Range("R10").Value = 1.847705673092716
$ws.Range("S10").Value = 2.020749367497032
$ws.Range("K14").Value = 16.86342592592595
$ws.Range("R14").Value = 1.92665172779809
$ws.Range("S14").Value = 2.116885095206829
$ws.Range("K15").Value = 16.86342592592595
$ws.Range("R15").Value = 3.408530896850205
$ws.Range("S15").Value = 4.236746585735966
$ws.Range("K16").Value = 16.86342592592595
$ws.Range("K26").Value = 12.67039049919483
$ws.Range("R26").Value = 1.847705673092716
$ws.Range("S26").Value = 2.020749367497032
$ws.Range("K27").Value = 12.67039049919483
$ws.Range("R27").Value = 3.135414186445925
$ws.Range("S27").Value = 3.817068050129936
$ws.Range("K28").Value = 12.67039049919483
$ws.Range("K42").Value = -1.819444444444444
$ws.Range("R42").Value = 2.455497817501559
$ws.Range("S42").Value = 2.843656807626497
$ws.Range("K43").Value = -1.819444444444444
$ws.Range("K44").Value = -1.819444444444444
$ws.Range("R44").Value = 1.618523362263702
$ws.Range("S44").Value = 1.746638928617865
$ws.Range("K45").Value = 13.62268518518517
$ws.Range("K46").Value = 13.62268518518517
$ws.Range("R46").Value = 1.865062221714807
$ws.Range("S46").Value = 2.041808739708676
$ws.Range("K47").Value = 13.62268518518517
$ws.Range("R47").Value = 3.193530074341067
$ws.Range("S47").Value = 3.904917482517481
$ws.Range("K51").Value = 13.0158303464755
$ws.Range("R51").Value = 3.15624931769735
$ws.Range("S51").Value = 3.848474349579903
$ws.Range("K52").Value = 13.0158303464755
$ws.Range("K53").Value = 13.0158303464755
$ws.Range("R53").Value = 1.853964204859962
$ws.Range("S53").Value = 2.02833814451736
$ws.Range("K54").Value = -1.819444444444444
$ws.Range("K55").Value = -1.819444444444444
$ws.Range("R55").Value = 2.455497817501559
$ws.Range("S55").Value = 2.843656807626497
$ws.Range("K56").Value = -1.819444444444444
$ws.Range("R56").Value = 1.618523362263702
$ws.Range("S56").Value = 1.746638928617865
$ws.Range("K60").Value = 13.0158303464755
$ws.Range("K61").Value = 13.0158303464755
$ws.Range("R61").Value = 1.853964204859962
$ws.Range("S61").Value = 2.02833814451736
$ws.Range("K62").Value = 13.0158303464755
$ws.Range("R62").Value = 3.15624931769735
$ws.Range("S62").Value = 3.848474349579903
$ws.Range("K66").Value = 19.36574074074073
$ws.Range("R66").Value = 3.595434716445165
$ws.Range("S66").Value = 4.534260101338181
$ws.Range("K67").Value = 19.36574074074073
$ws.Range("R67").Value = 1.977063465169192
$ws.Range("S67").Value = 2.178742498783586
$ws.Range("K68").Value = 19.36574074074073
$ws.Range("K69").Value = 13.62268518518517
$ws.Range("R69").Value = 3.193530074341067
$ws.Range("S69").Value = 3.904917482517481
$ws.Range("K70").Value = 13.62268518518517
$ws.Range("R70").Value = 1.865062221714807
$ws.Range("S70").Value = 2.041808739708676
$ws.Range("K71").Value = 13.62268518518517
$ws.Range("K75").Value = 19.65277777777778
$ws.Range("K76").Value = 19.65277777777778
$ws.Range("R76").Value = 1.983015294974508
$ws.Range("S76").Value = 2.18606997558991
$ws.Range("K77").Value = 19.65277777777778
$ws.Range("R77").Value = 3.618192955589586
$ws.Range("S77").Value = 4.571080550098231
$ws.Range("K93").Value = 13.0158303464755
$ws.Range("K94").Value = 13.0158303464755
$ws.Range("R94").Value = 3.15624931769735
$ws.Range("S94").Value = 3.848474349579903
$ws.Range("K95").Value = 13.0158303464755
$ws.Range("R95").Value = 1.853964204859962
$ws.Range("S95").Value = 2.02833814451736
$ws.Range("K99").Value = -1.819444444444444
$ws.Range("R99").Value = 2.455497817501559
$ws.Range("S99").Value = 2.843656807626497
$ws.Range("K100").Value = -1.819444444444444
$ws.Range("R100").Value = 1.618523362263702
$ws.Range("S100").Value = 1.746638928617865
$ws.Range("K101").Value = -1.819444444444444
$ws.Range("K105").Value = 5.486111111111112
$ws.Range("R105").Value = 2.756919486581097
$ws.Range("S105").Value = 3.263225806451613
$ws.Range("K106").Value = 5.486111111111112
$ws.Range("R106").Value = 1.726493341788205
$ws.Range("S106").Value = 1.874863921842289
$ws.Range("K107").Value = 5.486111111111112
